$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 (Feria Lagunitas de Puerto Montt -
# Espárragos weekly data). This pushes the existing rows 5-14 down to 6-15.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44481
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 300000000
$ws.Cells.Item(5, 7).Value = "Espárragos"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 1700
$ws.Cells.Item(5, 12).Value = 2000
$ws.Cells.Item(5, 13).Value = 1850
$ws.Cells.Item(5, 14).Value = '$/kilo'
$ws.Cells.Item(5, 15).Value = "Provincia de Linares"
$ws.Cells.Item(5, 16).Value = 1850
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
